$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 12.17913027012948
$ws.Range("C2").Value = 8.549402026451745
$ws.Range("D2").Value = 5.971302686218404
$ws.Range("E2").Value = 11.28459210549862
$ws.Range("G2").Value = 3.612409147005827
$ws.Range("I2").Value = 18.9814550723197
$ws.Range("M2").Value = 14.75491196416102
$ws.Range("N2").Value = 16.74192519004302
$ws.Range("O2").Value = 20.81998707484196

$ws.Range("B3").Value = 11.62472746345773
$ws.Range("C3").Value = 8.019774952927095
$ws.Range("D3").Value = 5.85021739935645
$ws.Range("E3").Value = 11.17211571382977
$ws.Range("G3").Value = 3.6148441378543
$ws.Range("I3").Value = 19.05159396053701
$ws.Range("M3").Value = 14.48125101321697
$ws.Range("N3").Value = 16.80370574409741
$ws.Range("O3").Value = 20.80448057940708

$ws.Range("B4").Value = 11.27206528691165
$ws.Range("C4").Value = 7.674788078515999
$ws.Range("D4").Value = 5.776417067218845
$ws.Range("E4").Value = 11.10678223467012
$ws.Range("G4").Value = 3.616418225123737
$ws.Range("I4").Value = 19.09973759178861
$ws.Range("M4").Value = 14.31424661633843
$ws.Range("N4").Value = 16.84345810946379
$ws.Range("O4").Value = 20.80142608029618

$ws.Range("B5").Value = 11.12547174578156
$ws.Range("C5").Value = 7.529234748234598
$ws.Range("D5").Value = 5.746527002748504
$ws.Range("E5").Value = 11.08112375437205
$ws.Range("G5").Value = 3.617079608688937
$ws.Range("I5").Value = 19.12062783140128
$ws.Range("M5").Value = 14.24654419253858
$ws.Range("N5").Value = 16.8601162552479
$ws.Range("O5").Value = 20.80180621281921

$ws.Range("B6").Value = 11.10096244955249
$ws.Range("C6").Value = 7.504766016667504
$ws.Range("D6").Value = 5.74157627860634
$ws.Range("E6").Value = 11.07692226813733
$ws.Range("G6").Value = 3.617190636626118
$ws.Range("I6").Value = 19.12417325786022
$ws.Range("M6").Value = 14.23532639724945
$ws.Range("N6").Value = 16.86291007445221
$ws.Range("O6").Value = 20.80196740752542

$ws.Range("B7").Value = 11.27009965265433
$ws.Range("C7").Value = 7.672845196616937
$ws.Range("D7").Value = 5.776013151538426
$ws.Range("E7").Value = 11.10643225102418
$ws.Range("G7").Value = 3.616427063988505
$ws.Range("I7").Value = 19.10001418489699
$ws.Range("M7").Value = 14.31333200078362
$ws.Range("N7").Value = 16.84368090784568
$ws.Range("O7").Value = 20.80142463061373

$ws.Range("B8").Value = 11.99062807341975
$ws.Range("C8").Value = 8.370898548393114
$ws.Range("D8").Value = 5.929466913643878
$ws.Range("E8").Value = 11.24505421755827
$ws.Range("G8").Value = 3.613232377108807
$ws.Range("I8").Value = 19.00458167921752
$ws.Range("M8").Value = 14.66039810276333
$ws.Range("N8").Value = 16.76285044640277
$ws.Range("O8").Value = 20.81329801680267

$ws.Range("B9").Value = 13.29886446893361
$ws.Range("C9").Value = 9.582583393817574
$ws.Range("D9").Value = 6.23270913034546
$ws.Range("E9").Value = 11.54511827197953
$ws.Range("G9").Value = 3.607591333798028
$ws.Range("I9").Value = 18.85797592837308
$ws.Range("M9").Value = 15.34498340328119
$ws.Range("N9").Value = 16.61871103406908
$ws.Range("O9").Value = 20.88787935504501

$ws.Range("B10").Value = 14.18788508172133
$ws.Range("C10").Value = 10.3770660748654
$ws.Range("D10").Value = 6.4542992592535
$ws.Range("E10").Value = 11.78083909517387
$ws.Range("G10").Value = 3.603822837291236
$ws.Range("I10").Value = 18.7753033890737
$ws.Range("M10").Value = 15.84498242147886
$ws.Range("N10").Value = 16.52148268588167
$ws.Range("O10").Value = 20.97383652563388

$ws.Range("B11").Value = 14.57532512062296
$ws.Range("C11").Value = 10.71778226633998
$ws.Range("D11").Value = 6.554324474216708
$ws.Range("E11").Value = 11.89094435657982
$ws.Range("G11").Value = 3.602189182794537
$ws.Range("I11").Value = 18.74319685186347
$ws.Range("M11").Value = 16.07071646558191
$ws.Range("N11").Value = 16.47911501344635
$ws.Range("O11").Value = 21.01965162762892

$ws.Range("B12").Value = 14.71950883161469
$ws.Range("C12").Value = 10.84383407271577
$ws.Range("D12").Value = 6.592049063697043
$ws.Range("E12").Value = 11.93301305999457
$ws.Range("G12").Value = 3.601582089092423
$ws.Range("I12").Value = 18.73183522690142
$ws.Range("M12").Value = 16.15586373508319
$ws.Range("N12").Value = 16.46333785106615
$ws.Range("O12").Value = 21.03795865374698

$ws.Range("B13").Value = 14.68857014784499
$ws.Range("C13").Value = 10.81681853730083
$ws.Range("D13").Value = 6.58393182193604
$ws.Range("E13").Value = 11.92393679669711
$ws.Range("G13").Value = 3.601712325510062
$ws.Range("I13").Value = 18.7342466529324
$ws.Range("M13").Value = 16.13754189820908
$ws.Range("N13").Value = 16.46672390898208
$ws.Range("O13").Value = 21.03397344274036

$ws.Range("B14").Value = 14.5872383962716
$ws.Range("C14").Value = 10.72821210658976
$ws.Range("D14").Value = 6.557431387659072
$ws.Range("E14").Value = 11.89439810005489
$ws.Range("G14").Value = 3.60213900598974
$ws.Range("I14").Value = 18.74224613305718
$ws.Range("M14").Value = 16.0777287413432
$ws.Range("N14").Value = 16.47781168141556
$ws.Range("O14").Value = 21.02113860922939

$ws.Range("B15").Value = 14.52483769682658
$ws.Range("C15").Value = 10.67355153753399
$ws.Range("D15").Value = 6.54117803126576
$ws.Range("E15").Value = 11.8763523604301
$ws.Range("G15").Value = 3.602401860570788
$ws.Range("I15").Value = 18.74724991700038
$ws.Range("M15").Value = 16.04104556269462
$ws.Range("N15").Value = 16.48463794041568
$ws.Range("O15").Value = 21.01340139509436

$ws.Range("B16").Value = 14.16221526938499
$ws.Range("C16").Value = 10.35438348048706
$ws.Range("D16").Value = 6.44774317656398
$ws.Range("E16").Value = 11.77369797033421
$ws.Range("G16").Value = 3.60393121809203
$ws.Range("I16").Value = 18.77751279123488
$ws.Range("M16").Value = 15.83018812702617
$ws.Range("N16").Value = 16.52428887822901
$ws.Range("O16").Value = 20.97097687762188

$ws.Range("B17").Value = 13.9353404848453
$ws.Range("C17").Value = 10.15328697885808
$ws.Range("D17").Value = 6.390194932002448
$ws.Range("E17").Value = 11.71143144621591
$ws.Range("G17").Value = 3.604890042386557
$ws.Range("I17").Value = 18.79749101844929
$ws.Range("M17").Value = 15.70033115651883
$ws.Range("N17").Value = 16.54908948759856
$ws.Range("O17").Value = 20.94666530232147

$ws.Range("B18").Value = 13.80325420832352
$ws.Range("C18").Value = 10.03567291285753
$ws.Range("D18").Value = 6.357023857744906
$ws.Range("E18").Value = 11.67589035801863
$ws.Range("G18").Value = 3.605449127928923
$ws.Range("I18").Value = 18.80949950072346
$ws.Range("M18").Value = 15.62548405099766
$ws.Range("N18").Value = 16.56352947178484
$ws.Range("O18").Value = 20.93331430788659

$ws.Range("B19").Value = 13.75826129404248
$ws.Range("C19").Value = 9.99551595618704
$ws.Range("D19").Value = 6.3457818217695
$ws.Range("E19").Value = 11.66390482188206
$ws.Range("G19").Value = 3.60563973106375
$ws.Range("I19").Value = 18.8136540885889
$ws.Range("M19").Value = 15.60011787376029
$ws.Range("N19").Value = 16.56844875214099
$ws.Range("O19").Value = 20.92890270412274

$ws.Range("B20").Value = 13.95965735328551
$ws.Range("C20").Value = 10.17489568436428
$ws.Range("D20").Value = 6.396328677717979
$ws.Range("E20").Value = 11.7180318544369
$ws.Range("G20").Value = 3.604787188282722
$ws.Range("I20").Value = 18.7953107000734
$ws.Range("M20").Value = 15.7141715206775
$ws.Range("N20").Value = 16.54643128313689
$ws.Range("O20").Value = 20.94918791096615

$ws.Range("B21").Value = 14.61707133201243
$ws.Range("C21").Value = 10.75431851106509
$ws.Range("D21").Value = 6.565219669549663
$ws.Range("E21").Value = 11.90306448932083
$ws.Range("G21").Value = 3.602013367004141
$ws.Range("I21").Value = 18.73987483598486
$ws.Range("M21").Value = 16.09530702722369
$ws.Range("N21").Value = 16.47454770907496
$ws.Range("O21").Value = 21.02488258236489

$ws.Range("B22").Value = 15.03194269697734
$ws.Range("C22").Value = 11.11569716641827
$ws.Range("D22").Value = 6.674691774279186
$ws.Range("E22").Value = 12.02615685938172
$ws.Range("G22").Value = 3.60026772541887
$ws.Range("I22").Value = 18.70828876141827
$ws.Range("M22").Value = 16.34242143914887
$ws.Range("N22").Value = 16.42912078674003
$ws.Range("O22").Value = 21.07993192679836

$ws.Range("B23").Value = 14.81189693813038
$ws.Range("C23").Value = 10.92440377798432
$ws.Range("D23").Value = 6.616360469046007
$ws.Range("E23").Value = 11.96027529529212
$ws.Range("G23").Value = 3.601193277777277
$ws.Range("I23").Value = 18.7247201820733
$ws.Range("M23").Value = 16.21074032188417
$ws.Range("N23").Value = 16.45322426730105
$ws.Range("O23").Value = 21.05004347758032

$ws.Range("B24").Value = 13.94866883863509
$ws.Range("C24").Value = 10.16513262187219
$ws.Range("D24").Value = 6.3935558754778
$ws.Range("E24").Value = 11.71504700712011
$ws.Range("G24").Value = 3.604833664182211
$ws.Range("I24").Value = 18.79629479406124
$ws.Range("M24").Value = 15.70791487910137
$ws.Range("N24").Value = 16.54763249099709
$ws.Range("O24").Value = 20.94804548886055

$ws.Range("B25").Value = 12.95710541903123
$ws.Range("C25").Value = 9.2716216948432
$ws.Range("D25").Value = 6.150701823601538
$ws.Range("E25").Value = 11.46111534956402
$ws.Range("G25").Value = 3.60905105337493
$ws.Range("I25").Value = 18.89326254416778
$ws.Range("M25").Value = 15.15992575336992
$ws.Range("N25").Value = 16.65617552304449
$ws.Range("O25").Value = 20.86221659991494
